# Add the Emon, CF3hr & Esubhr reffclws combinations of the PEXTRA
# "reffclws" variable to the missing-identified-variables list (#564).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$variable   = "reffclws"
$commentTxt = "From the AerChem side there is interest in Emon reffclws. This variable is identified as the already available PEXTRA variable with the table 126 grib code 126021, i.e. proposing to add reffclws as 21.126 to ifspar.json.  Note that this variable is not requested by CMIP6 AerChem, and that reffclws not occurs in any CMIP6 data request of the experiments in which any EC-Earth3* configuration participates. See further #564."
$author     = "Thomas"

$tables = @("Emon", "CF3hr", "Esubhr")

$row = 85
foreach ($table in $tables) {
    $ws.Range("A$row").Value = $table

    $ws.Range("B$row").Value = $variable
    $ws.Range("B$row").WrapText = $true

    $ws.Range("H$row").Value = $commentTxt
    $ws.Range("I$row").Value = $author

    $row = $row + 1
}

$ws.Range("A85:I87").Select()
